$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 11
$ws.Range("D2").Value = 27
$ws.Range("F2").Value = 1065.864898671727
$ws.Range("G2").Value = -1245.957256848989
$ws.Range("H2").Value = 0.1061583337523293
$ws.Range("I2").Value = -0.1240457479771876
$ws.Range("J2").Value = -0.05740772431864322
$ws.Range("K2").Value = 11724.513885389
$ws.Range("L2").Value = -33640.84593492271
$ws.Range("M2").Value = 0.2894736842105263
$ws.Range("N2").Value = 0.8554586385790522
$ws.Range("O2").Value = 0.3485201860877621
$ws.Range("P2").Value = -21916.33204953371
$ws.Range("C3").Value = 110
$ws.Range("D3").Value = 125
$ws.Range("E3").Value = 235
$ws.Range("F3").Value = 1940.796189685295
$ws.Range("G3").Value = -932.3577508378729
$ws.Range("H3").Value = 0.1940126954695114
$ws.Range("I3").Value = -0.09331122958962974
$ws.Range("J3").Value = 0.04118082043805334
$ws.Range("K3").Value = 213487.5808653825
$ws.Range("L3").Value = -116544.7188547342
$ws.Range("M3").Value = 0.4680851063829787
$ws.Range("N3").Value = 2.081600316982594
$ws.Range("O3").Value = 1.831808278944683
$ws.Range("P3").Value = 96942.86201064831
$ws.Range("C4").Value = 504
$ws.Range("D4").Value = 461
$ws.Range("E4").Value = 965
$ws.Range("F4").Value = 2627.3161530074
$ws.Range("G4").Value = -993.0214025778253
$ws.Range("H4").Value = 0.2710758953053513
$ws.Range("I4").Value = -0.09831450929959681
$ws.Range("J4").Value = 0.09461063465987872
$ws.Range("K4").Value = 1324167.341115729
$ws.Range("L4").Value = -457782.8665883777
$ws.Range("M4").Value = 0.522279792746114
$ws.Range("N4").Value = 2.645779986400133
$ws.Range("O4").Value = 2.892566405955891
$ws.Range("P4").Value = 866384.4745273512
$ws.Range("C5").Value = 194
$ws.Range("D5").Value = 336
$ws.Range("E5").Value = 530
$ws.Range("F5").Value = 2140.966170232832
$ws.Range("G5").Value = -970.0211297631088
$ws.Range("H5").Value = 0.2122189173511524
$ws.Range("I5").Value = -0.09520616080491509
$ws.Range("J5").Value = 0.0173230187465512
$ws.Range("K5").Value = 415347.4370251691
$ws.Range("L5").Value = -325927.0996004042
$ws.Range("M5").Value = 0.3660377358490566
$ws.Range("N5").Value = 2.207133540230905
$ws.Range("O5").Value = 1.274356865490463
$ws.Range("P5").Value = 89420.33742476493
$ws.Range("C6").Value = 147
$ws.Range("D6").Value = 266
$ws.Range("E6").Value = 413
$ws.Range("F6").Value = 2046.452982136613
$ws.Range("G6").Value = -846.7055717687756
$ws.Range("H6").Value = 0.2043651906374828
$ws.Range("I6").Value = -0.08451588313545165
$ws.Range("J6").Value = 0.01830619397016903
$ws.Range("K6").Value = 300828.588374082
$ws.Range("L6").Value = -225223.6820904944
$ws.Range("M6").Value = 0.3559322033898305
$ws.Range("N6").Value = 2.416959389863886
$ws.Range("O6").Value = 1.335688083872147
$ws.Range("P6").Value = 75604.90628358762
$ws.Range("C7").Value = 141
$ws.Range("D7").Value = 194
$ws.Range("E7").Value = 335
$ws.Range("F7").Value = 2221.640206844884
$ws.Range("G7").Value = -818.5710704103749
$ws.Range("H7").Value = 0.2298554816768995
$ws.Range("I7").Value = -0.08202277077029688
$ws.Range("J7").Value = 0.04924538921494098
$ws.Range("K7").Value = 313251.2691651287
$ws.Range("L7").Value = -158802.7876596127
$ws.Range("M7").Value = 0.4208955223880597
$ws.Range("N7").Value = 2.714046815423257
$ws.Range("O7").Value = 1.972580417395254
$ws.Range("P7").Value = 154448.481505516
$ws.Range("C8").Value = 132
$ws.Range("D8").Value = 255
$ws.Range("E8").Value = 387
$ws.Range("F8").Value = 2185.868283077339
$ws.Range("G8").Value = -834.1021821127189
$ws.Range("H8").Value = 0.2187874786231606
$ws.Range("I8").Value = -0.08304863756446565
$ws.Range("J8").Value = 0.01990321601891077
$ws.Range("K8").Value = 288534.6133662087
$ws.Range("L8").Value = -212696.0564387435
$ws.Range("M8").Value = 0.3410852713178295
$ws.Range("N8").Value = 2.620624103321127
$ws.Range("O8").Value = 1.356558359366229
$ws.Range("P8").Value = 75838.55692746525
$ws.Range("C9").Value = 18
$ws.Range("D9").Value = 24
$ws.Range("F9").Value = 1294.276291301255
$ws.Range("G9").Value = -756.0723005579936
$ws.Range("H9").Value = 0.1258336468436882
$ws.Range("I9").Value = -0.07714004066264141
$ws.Range("J9").Value = 0.00984868255435703
$ws.Range("K9").Value = 23296.9732434226
$ws.Range("L9").Value = -18145.73521339185
$ws.Range("M9").Value = 0.4285714285714285
$ws.Range("N9").Value = 1.711841963190635
$ws.Range("O9").Value = 1.283881472392976
$ws.Range("P9").Value = 5151.238030030749
$ws.Range("C10").Value = 397
$ws.Range("D10").Value = 532
$ws.Range("E10").Value = 929
$ws.Range("F10").Value = 2754.67662934896
$ws.Range("G10").Value = -1002.523629543307
$ws.Range("H10").Value = 0.2748679142881063
$ws.Range("I10").Value = -0.09941061602228903
$ws.Range("J10").Value = 0.06053403040744949
$ws.Range("K10").Value = 1093606.621851538
$ws.Range("L10").Value = -533342.5709170396
$ws.Range("M10").Value = 0.4273412271259419
$ws.Range("N10").Value = 2.747742345588237
$ws.Range("O10").Value = 2.050476900749118
$ws.Range("P10").Value = 560264.0509344984
$ws.Range("C11").Value = 255
$ws.Range("D11").Value = 330
$ws.Range("E11").Value = 585
$ws.Range("F11").Value = 1976.227219895312
$ws.Range("G11").Value = -907.2027938939198
$ws.Range("H11").Value = 0.1984282379405465
$ws.Range("I11").Value = -0.09041322838037601
$ws.Range("J11").Value = 0.03549202616976976
$ws.Range("K11").Value = 503937.9410733046
$ws.Range("L11").Value = -299376.9219849936
$ws.Range("M11").Value = 0.4358974358974359
$ws.Range("N11").Value = 2.178374265596005
$ws.Range("O11").Value = 1.683289205233277
$ws.Range("P11").Value = 204561.019088311
$ws.Range("C12").Value = 161
$ws.Range("D12").Value = 231
$ws.Range("E12").Value = 392
$ws.Range("F12").Value = 2233.501933587452
$ws.Range("G12").Value = -925.1404744456071
$ws.Range("H12").Value = 0.226940156990908
$ws.Range("I12").Value = -0.09221070856189445
$ws.Range("J12").Value = 0.03886911121872084
$ws.Range("K12").Value = 359593.8113075797
$ws.Range("L12").Value = -213707.4495969353
$ws.Range("M12").Value = 0.4107142857142857
$ws.Range("N12").Value = 2.414230049686113
$ws.Range("O12").Value = 1.682645186144866
$ws.Range("P12").Value = 145886.3617106444
$ws.Range("C13").Value = 271
$ws.Range("D13").Value = 438
$ws.Range("E13").Value = 709
$ws.Range("F13").Value = 1935.847265955373
$ws.Range("G13").Value = -925.9098326102682
$ws.Range("H13").Value = 0.1977033107758932
$ws.Range("I13").Value = -0.09219071664821937
$ws.Range("J13").Value = 0.01861503995535538
$ws.Range("K13").Value = 524614.6090739062
$ws.Range("L13").Value = -405548.5066832977
$ws.Range("M13").Value = 0.382228490832158
$ws.Range("N13").Value = 2.090751386123583
$ws.Range("O13").Value = 1.293592752601577
$ws.Range("P13").Value = 119066.1023906085
$ws.Range("C14").Value = 507
$ws.Range("D14").Value = 568
$ws.Range("E14").Value = 1075
$ws.Range("F14").Value = 2066.836578176029
$ws.Range("G14").Value = -819.1683359369599
$ws.Range("H14").Value = 0.2075064045666483
$ws.Range("I14").Value = -0.08199896206661708
$ws.Range("J14").Value = 0.05453984805716494
$ws.Range("K14").Value = 1047886.145135246
$ws.Range("L14").Value = -465287.6148121934
$ws.Range("M14").Value = 0.4716279069767442
$ws.Range("N14").Value = 2.523091393433797
$ws.Range("O14").Value = 2.252125592378405
$ws.Range("P14").Value = 582598.5303230529
$ws.Range("C15").Value = 163
$ws.Range("D15").Value = 352
$ws.Range("E15").Value = 515
$ws.Range("F15").Value = 2258.922078427619
$ws.Range("G15").Value = -966.2108762869353
$ws.Range("H15").Value = 0.2246646777007115
$ws.Range("I15").Value = -0.09646070284454954
$ws.Range("J15").Value = 0.00517703895909613
$ws.Range("K15").Value = 368204.2987837017
$ws.Range("L15").Value = -340106.2284530012
$ws.Range("M15").Value = 0.3165048543689321
$ws.Range("N15").Value = 2.337918288716083
$ws.Range("O15").Value = 1.082615571195231
$ws.Range("P15").Value = 28098.07033070043
$ws.Range("C16").Value = 134
$ws.Range("D16").Value = 259
$ws.Range("E16").Value = 393
$ws.Range("F16").Value = 1811.515751749597
$ws.Range("G16").Value = -981.2746221879318
$ws.Range("H16").Value = 0.1817883529152428
$ws.Range("I16").Value = -0.0975941553730008
$ws.Range("J16").Value = -0.002333961707289194
$ws.Range("K16").Value = 242743.110734446
$ws.Range("L16").Value = -254150.1271466743
$ws.Range("M16").Value = 0.3409669211195929
$ws.Range("N16").Value = 1.846084379223514
$ws.Range("O16").Value = 0.9551170147334016
$ws.Range("P16").Value = -11407.01641222829
$ws.Range("C17").Value = 439
$ws.Range("D17").Value = 523
$ws.Range("E17").Value = 962
$ws.Range("F17").Value = 1648.779263303026
$ws.Range("G17").Value = -929.278701240899
$ws.Range("H17").Value = 0.1651016915794723
$ws.Range("I17").Value = -0.09211484013429698
$ws.Range("J17").Value = 0.02526359793466836
$ws.Range("K17").Value = 723814.0965900293
$ws.Range("L17").Value = -486012.7607489906
$ws.Range("M17").Value = 0.4563409563409563
$ws.Range("N17").Value = 1.774257024401132
$ws.Range("O17").Value = 1.489290313025043
$ws.Range("P17").Value = 237801.3358410387
$ws.Range("C18").Value = 392
$ws.Range("D18").Value = 643
$ws.Range("E18").Value = 1035
$ws.Range("F18").Value = 2002.039561261423
$ws.Range("G18").Value = -921.6669757236119
$ws.Range("H18").Value = 0.1995506913983129
$ws.Range("I18").Value = -0.09182953350317111
$ws.Range("J18").Value = 0.01852896713584513
$ws.Range("K18").Value = 784799.508014478
$ws.Range("L18").Value = -592631.8653902824
$ws.Range("M18").Value = 0.378743961352657
$ws.Range("N18").Value = 2.172194093956331
$ws.Range("O18").Value = 1.324261407201994
$ws.Range("P18").Value = 192167.6426241957
$ws.Range("C19").Value = 323
$ws.Range("D19").Value = 613
$ws.Range("E19").Value = 936
$ws.Range("F19").Value = 2101.4519154231
$ws.Range("G19").Value = -923.8385216330095
$ws.Range("H19").Value = 0.2071412624601392
$ws.Range("I19").Value = -0.09215566831216562
$ws.Range("J19").Value = 0.01112735373853354
$ws.Range("K19").Value = 678768.9686816609
$ws.Range("L19").Value = -566313.0137610348
$ws.Range("M19").Value = 0.3450854700854701
$ws.Range("N19").Value = 2.274696135974607
$ws.Range("O19").Value = 1.198575614877321
$ws.Range("P19").Value = 112455.9549206261
$ws.Range("C20").Value = 152
$ws.Range("E20").Value = 478
$ws.Range("F20").Value = 2302.774924639733
$ws.Range("G20").Value = -1050.799571582465
$ws.Range("H20").Value = 0.2264910358078325
$ws.Range("I20").Value = -0.1043882037536548
$ws.Range("J20").Value = 0.0008286255629687808
$ws.Range("K20").Value = 350021.7885452395
$ws.Range("L20").Value = -342560.6603358834
$ws.Range("M20").Value = 0.3179916317991632
$ws.Range("N20").Value = 2.191450193657618
$ws.Range("O20").Value = 1.021780458392509
$ws.Range("P20").Value = 7461.128209356102
$ws.Range("C21").Value = 394
$ws.Range("D21").Value = 905
$ws.Range("E21").Value = 1299
$ws.Range("F21").Value = 2283.165363628508
$ws.Range("G21").Value = -1214.332938672681
$ws.Range("H21").Value = 0.2283173661189831
$ws.Range("I21").Value = -0.116914055953673
$ws.Range("J21").Value = -0.01220183093702418
$ws.Range("K21").Value = 899567.1532696316
$ws.Range("L21").Value = -1098971.309498776
$ws.Range("M21").Value = 0.3033102386451116
$ws.Range("N21").Value = 1.88018070738006
$ws.Range("O21").Value = 0.8185538107267882
$ws.Range("P21").Value = -199404.1562291442
